# UC007 workbook: swap the "step" content between TC3 and TC4 test cases.
#
# Before:
#   TC3 step 2: "Chefe Seleciona um usuário para filtrar as autorizações de
#                pagamento associadas a ele; e Submete a busca ao sistema."
#               / "SYSTEM Filtra os registros (autorizações de pagamento
#                pendentes) e exibe apenas aqueles atribuídos ao usuário
#                selecionado."
#   TC4 step 2: "Chefe Clica para realizar a autorização de pagamento."
#               / "SYSTEM Apresenta a tela de Registrar Autorizações de
#                Pagamento"
#
# After (this edit):
#   TC3 step 2: "Chefe Clica para realizar a autorização de pagamento."
#               / "SYSTEM Apresenta a tela de Registrar Autorizações de
#                Pagamento"
#   TC4 step 2: "Chefe Seleciona um usuário para filtrar as autorizações de
#                pagamento associadas a ele; e Submete a busca ao sistema."
#               / "SYSTEM Filtra os registros (autorizações de pagamento
#                pendentes) e exibe apenas aqueles atribuídos ao usuário
#                selecionado."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$tc3StepText        = "Chefe Seleciona um usuário para filtrar as autorizações de pagamento associadas a ele; e Submete a busca ao sistema."
$tc3ExpectedText    = "SYSTEM Filtra os registros (autorizações de pagamento pendentes) e exibe apenas aqueles atribuídos ao usuário selecionado."
$tc4StepText        = "Chefe Clica para realizar a autorização de pagamento."
$tc4ExpectedText    = "SYSTEM Apresenta a tela de Registrar Autorizações de Pagamento"

# Row 28 holds TC3's second step (B28 = step, D28 = expected result).
# Row 36 holds TC4's second step (B36 = step, D36 = expected result).
# Swap the contents so TC3 gets the "Clica para realizar..." step and
# TC4 gets the "Seleciona um usuário..." step.

$ws.Range("B28").Value = $tc4StepText
$ws.Range("D28").Value = $tc4ExpectedText

$ws.Range("B36").Value = $tc3StepText
$ws.Range("D36").Value = $tc3ExpectedText

$wb.Save()
